# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 4.33
$ws.Range("K3").Value = 1.8
$ws.Range("U3").Value = 5.6
$ws.Range("AS3").Value = 51

# Row 5
$ws.Range("O5").Value = 1.73
$ws.Range("P5").Value = 2
$ws.Range("S5").Value = 3.5
$ws.Range("T5").Value = 1.3
$ws.Range("U5").Value = 6.6

# Row 6
$ws.Range("N6").Value = 5
$ws.Range("Z6").Value = 2.08

# Row 8
$ws.Range("G8").Value = 2.6
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.4
$ws.Range("L8").Value = 3.4
$ws.Range("N8").Value = 8.5
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("S8").Value = 2.15
$ws.Range("T8").Value = 1.67
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 12
$ws.Range("AE8").Value = 10
$ws.Range("AF8").Value = 26
$ws.Range("AI8").Value = 8.5
$ws.Range("AM8").Value = 301
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 11
$ws.Range("AQ8").Value = 26
$ws.Range("AR8").Value = 23
$ws.Range("AS8").Value = 34

# Row 10
$ws.Range("G10").Value = 2.12
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.15
$ws.Range("J10").Value = 2.75
$ws.Range("K10").Value = 2.12
$ws.Range("L10").Value = 3.65
$ws.Range("N10").Value = 6.9
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.05
$ws.Range("S10").Value = 1.98
$ws.Range("W10").Value = 3.35
$ws.Range("X10").Value = 1.28
$ws.Range("Y10").Value = 1.4
$ws.Range("Z10").Value = 2.75
$ws.Range("AC10").Value = 7.1
$ws.Range("AD10").Value = 10
$ws.Range("AE10").Value = 8.75
$ws.Range("AF10").Value = 20
$ws.Range("AG10").Value = 18
$ws.Range("AH10").Value = 30
$ws.Range("AI10").Value = 6.9
$ws.Range("AJ10").Value = 6.3
$ws.Range("AN10").Value = 9.5
$ws.Range("AO10").Value = 17
$ws.Range("AP10").Value = 11.25
$ws.Range("AQ10").Value = 40
$ws.Range("AR10").Value = 28
$ws.Range("AS10").Value = 37

# Row 11
$ws.Range("G11").Value = 2.05
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.63
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 1.03
$ws.Range("O11").Value = 1.22
$ws.Range("S11").Value = 1.82
$ws.Range("T11").Value = 1.92
$ws.Range("X11").Value = 1.3
$ws.Range("Y11").Value = 1.36
$ws.Range("Z11").Value = 3
$ws.Range("AC11").Value = 8
$ws.Range("AD11").Value = 10
$ws.Range("AF11").Value = 19
$ws.Range("AO11").Value = 19
$ws.Range("AQ11").Value = 41
$ws.Range("AR11").Value = 29

# Row 12
$ws.Range("M12").Value = 1.08
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.37
$ws.Range("Q12").Value = 1.94
$ws.Range("R12").Value = 1.79
$ws.Range("X12").Value = 1.11

# Row 13
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 2.6
$ws.Range("K13").Value = 2.05
$ws.Range("L13").Value = 3.4
$ws.Range("M13").Value = 1.04
$ws.Range("O13").Value = 1.3
$ws.Range("X13").Value = 1.22
$ws.Range("AN13").Value = 8.5
$ws.Range("AO13").Value = 13
$ws.Range("AR13").Value = 23
